$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price / volume data (values + a couple of row re-orderings)

# Cells whose new text looks like a plain number (e.g. "1.00", "0.999") must be
# forced to Text format first so Excel does not silently convert them to a numeric value
# and drop formatting such as trailing zeros.
$textForcedCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D14",
    "D18",
    "D19",
    "D21",
    "D24",
    "D28",
    "D30",
    "D32",
    "D33",
    "D34",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$cellValues = @{
    "D2" = "58.646.74"
    "E2" = "  +2.28%  "
    "D3" = "3.106.27"
    "E3" = "  +0.79%  "
    "D4" = "1.00"
    "E4" = "  +0.00%  "
    "D5" = "527.98"
    "E5" = "  +2.29%  "
    "D6" = "143.79"
    "E6" = "  +1.58%  "
    "D7" = "1.00"
    "E7" = "  +0.02%  "
    "E8" = "  +1.82%  "
    "E9" = "  +1.26%  "
    "E10" = "  +0.52%  "
    "E11" = "  +3.01%  "
    "D12" = "3.638.23"
    "E12" = "  +0.78%  "
    "E13" = "  +1.02%  "
    "D14" = "26.95"
    "E14" = "  +5.22%  "
    "E15" = "  +1.81%  "
    "D16" = "58.669.72"
    "E16" = "  +2.13%  "
    "B17" = "WrappedEther"
    "C17" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D17" = "3.097.79"
    "E17" = "  +0.51%  "
    "B18" = "Polkadot"
    "C18" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D18" = "6.14"
    "E18" = "  +0.37%  "
    "D19" = "12.94"
    "E19" = "  -1.31%  "
    "E20" = "  -0.34%  "
    "D21" = "341.64"
    "E21" = "  +2.28%  "
    "E22" = "  +0.00%  "
    "E23" = "  +1.01%  "
    "D24" = "65.98"
    "E24" = "  +0.08%  "
    "E25" = "  +0.43%  "
    "E26" = "  -0.05%  "
    "E27" = "  +0.16%  "
    "D28" = "6.70"
    "E28" = "  +4.37%  "
    "E29" = "  +2.02%  "
    "D30" = "1.87"
    "E30" = "  +3.04%  "
    "E31" = "  +3.36%  "
    "D32" = "21.03"
    "E32" = "  +0.97%  "
    "D33" = "153.97"
    "E33" = "  -0.01%  "
    "D34" = "4.67"
    "E34" = "  +2.90%  "
    "E35" = "  +2.84%  "
    "D36" = "27.18"
    "E36" = "  -2.34%  "
    "D37" = "1.31"
    "E37" = "  +4.44%  "
    "D38" = "0.0679"
    "E38" = "  +0.37%  "
    "D39" = "3.145.86"
    "E39" = "  +0.80%  "
    "D40" = "3.91"
    "E40" = "  +1.60%  "
    "B41" = "Mantle"
    "C41" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "D41" = "0.677"
    "E41" = "  +0.88%  "
    "B42" = "OKB"
    "C42" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D42" = "36.83"
    "E42" = "  +0.16%  "
    "D43" = "1.50"
    "E43" = "  +7.50%  "
    "D44" = "0.999"
    "E44" = "  +0.03%  "
    "D45" = "2.295.14"
    "E45" = "  +0.05%  "
    "E46" = "  +1.43%  "
    "D47" = "21.01"
    "E47" = "  +4.55%  "
    "D48" = "0.973"
    "E48" = "  +3.24%  "
    "E49" = "  +1.82%  "
    "D50" = "0.752"
    "E50" = "  +9.35%  "
    "D51" = "268.69"
    "E51" = "  +8.05%  "
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}
